$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1, J1 ---
# Copy formatting (style) from the existing header cell H1, then set new values/text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data columns I2:J46 ---
$iValues = @(8,6,6,9,7,6,6,5,6,9,9,4,9,5,8,6,5,7,4,7,6,7,7,8,8,5,8,5,5,5,6,7,8,6,7,8,7,7,8,6,7,5,6,7,6)
$jValues = @(8,7,6,9,7,6,7,5,7,9,9,4,9,6,8,6,5,7,5,7,6,7,7,8,9,5,8,6,5,5,8,7,8,6,8,8,8,7,8,6,7,5,6,7,6)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}

Write-Output "applied I0/IF columns"
